$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the s1Protocol value (column H, data rows 2-27):
#    "E7760" -> "E7420"
$rng = $ws.Range("H2:H27")
$rng.Value = "E7420"

# 2. Clean up the cell formatting on that same range: force the font to an
#    explicit "automatic" (theme) colour instead of the previously explicit
#    black, producing a fresh cell style for H2:H27.
$rng.Font.Name = "Arial"
$rng.Font.Size = 10
$rng.Font.ThemeColor = 1

# 3. Reflect the new working selection (the sheet was left focused on the
#    s1Protocol column after the cleanup).
$ws.Range("H2:H27").Select() | Out-Null
